$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final dataset for rows 2-31 (A:H): 6 new rows were inserted at the top
# (pushing the former rows 2-21 down to rows 8-27, i.e. new rows 8-27 carry
# the previous C:H measurements), and 4 more new rows were appended at the
# bottom (rows 28-31). Column A (timestamp) and column B (label) follow the
# existing per-row pattern and are simply rewritten here together with the
# rest so the whole block ends up self consistent.
$rows = @(
    ,(0, "falling", -2.778316736221313, 6.619067668914795, -0.6493567824363708, -0.0410806871950626, -0.0163406450301408, -0.015118914656341)
    ,(100, "falling", -2.981808757781982, 6.65326156616211, -0.8340041637420658, -0.022754730656743, -0.0708603709936142, -0.0030543261673301)
    ,(200, "falling", -2.862793350219726, 6.595655870437622, -1.086881220340729, 0.0134390350431203, -0.0337503030896186, 0.0230601634830236)
    ,(300, "falling", -2.658387470245361, 6.523494625091553, -1.100642728805541, -0.0311541277915239, -0.0189368221908807, 0.0268780700862407)
    ,(400, "falling", -2.623731708526611, 6.502191925048828, -0.9057361066341401, -0.0083993971347808, -0.0177150927484035, 0.0058032199740409)
    ,(500, "falling", -2.912428379058838, 6.643356800079346, -0.9624049067497252, -0.0368046313524246, -0.0198531206697225, 0.0143553335219621)
    ,(600, "falling", -2.799852418899536, 6.562464380264283, -1.04518318772316, -0.0303905457258224, 0.0251981914043426, -0.0348193198442459)
    ,(700, "falling", -2.921730804443359, 6.461819076538085, -1.420121455192567, 0.0346666015684604, 0.0082466807216405, -0.0271835029125213)
    ,(800, "falling", -3.64505925178528, 6.319645357131958, -2.137768888473512, 0.1174388453364372, -0.140499010682106, 0.1302670091390609)
    ,(900, "falling", -3.560260105133056, 5.857164859771728, -2.574245262145996, -0.1805106848478317, -0.266184538602829, -0.0308486949652433)
    ,(1000, "falling", -3.488956928253174, 5.242366313934325, -3.541561126708986, -0.1626428663730621, -0.379652738571167, 0.0626136884093284)
    ,(1100, "falling", -3.976312160491945, 4.683113002777098, -4.901833915710452, 0.0490219369530677, -0.1159116774797439, 0.3532328307628631)
    ,(1200, "falling", -4.27849884033203, 4.203760409355163, -6.137033796310429, 0.0462730415165424, -0.1301143020391464, 0.0424551330506801)
    ,(1300, "falling", -4.059399700164795, 3.346987175941464, -6.092747974395745, 0.1779145002365112, 0.9054549932479858, -1.665982246398926)
    ,(1400, "falling", -4.056312322616568, 2.196766674518599, -0.6121338844299181, -0.7797694802284241, 2.544864654541016, -3.356246471405029)
    ,(1500, "falling", 1.137712478637695, 9.979324340820312, 6.995259761810303, -1.158047795295715, 0.633009135723114, -1.008691191673279)
    ,(1600, "falling", -5.27030158042908, 6.398296785354612, -1.064579951763156, -0.3301726579666137, -2.484847068786621, 0.2561052441596985)
    ,(1700, "falling", 9.112322330474868, 4.601176929473875, -2.501108241081238, 0.0949895456433296, 0.5517640113830566, 0.8726209998130798)
    ,(1800, "falling", 5.462616539001423, 5.74287872314454, -1.951083517074581, 0.1779145002365112, -0.2057088762521743, -0.0172569435089826)
    ,(1900, "falling", 3.208791923522958, 6.954292774200445, -1.887256383895876, -0.3599523603916168, -0.3149010241031647, -0.7423539757728577)
    ,(2000, "falling", 4.056319236755372, 5.56828045845031, -1.578303456306456, -0.022754730656743, 0.0948368310928344, -0.1800525337457656)
    ,(2100, "falling", 3.340787315368646, 4.649471092224126, -1.28086395263672, -0.0279470849782228, 0.1127046346664428, 0.2240348309278488)
    ,(2200, "falling", 1.86997218132019, 6.004907035827637, -1.588819026947022, -0.076052725315094, -0.1412625908851623, 0.1592831164598465)
    ,(2300, "falling", 1.657377624511718, 6.300898742675783, -1.739501070976257, -0.0526871271431446, -0.1530217379331588, -0.0560468845069408)
    ,(2400, "falling", 1.416164684295658, 6.649864673614499, -1.776257181167603, 0.0149661982432007, 0.0545197241008281, -0.081092357635498)
    ,(2500, "falling", 2.375874519348145, 6.059439659118652, -1.791717648506165, -0.0088575463742017, 0.0604756586253643, 0.0668897479772567)
    ,(2600, "falling", 1.882777214050291, 5.883625364303589, -1.771220207214356, -0.0070249503478407, 0.098502017557621, 0.1296561509370803)
    ,(2700, "falling", 2.030305290222169, 6.322917270660403, -1.870007395744325, -0.0024434609804302, -0.0125227374956011, -0.0360410511493682)
    ,(2800, "falling", 2.015394306182861, 6.204385328292845, -1.822170174121856, 0.016951510682702, 0.0704022198915481, 0.0198531206697225)
    ,(2900, "falling", 1.990867900848389, 6.087325572967529, -1.848094010353089, -0.00335975876078, -0.0123700210824608, 0.0325285755097866)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $row = $rows[$i]
    $targetRow = 2 + $i
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws.Cells.Item($targetRow, 1 + $j).Value = $row[$j]
    }
}
